$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.254.92'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = '2.306.57'
$ws.Range('E3').Value = '  -1.80%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('E7').Value = '  -1.42%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.610'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.83'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0908'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.51'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.17%  '
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.976'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.46'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.86%  '
$ws.Range('D16').Value = '2.655.49'
$ws.Range('E16').Value = '  -1.81%  '
$ws.Range('D17').Value = '2.297.57'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Value = '42.135.13'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.73'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.65'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.34%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '280.94'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.66%  '
$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.55'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.40%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.67%  '
$ws.Range('E26').Value = '  +0.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.39%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.34%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.27'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.94'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '165.57'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0886'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.92'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.91'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.120'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.07%  '
$ws.Range('E36').Value = '  +1.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.66'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.95'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0354'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.66'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '101.88'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +17.32%  '
$ws.Range('E42').Value = '  +1.90%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '71.01'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('E44').Value = '  -3.30%  '
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '117.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.16'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '78.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.34'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.72%  '
$ws.Range('E51').Value = '  +2.33%  '
